$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap itemname, FRDSKF (M), and MHKSKF (U) values between row 5 and row 6
$ws.Range("D5").Value = "Zithrox 20ml Powder for Suspension"
$ws.Range("D6").Value = "Zithrox 50ml Powder for Suspension"

$ws.Range("M5").Value = 3
$ws.Range("M6").Value = ""

$ws.Range("U5").Value = ""
$ws.Range("U6").Value = 344
